$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new row of form data (a new submission appended below the
# existing rows) that was added when "writing to spreadsheet on save"
# was implemented on mainPage.

$ws.Range("A4").Value = "asd"
$ws.Range("B4").Value = "asd"
$ws.Range("C4").Value = "asd"
$ws.Range("D4").Value = "asd"
$ws.Range("E4").Value = "asd"

# Force this one to be stored as text (not a number) since the source
# field is a free-text rate-limit entry.
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "120"

# The remaining columns of the submitted row were left blank by the
# user, but the row still has an entry (cell) for every column.
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Style = "Normal"
$ws.Range("I4").Style = "Normal"
$ws.Range("J4").Style = "Normal"
$ws.Range("K4").Style = "Normal"
$ws.Range("L4").Style = "Normal"
$ws.Range("M4").Style = "Normal"
